# Edit: add "metadata" sheet after "data", and refresh the F-column
# (time_taken) timestamps on "data" to reflect a later panel re-query.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- 1. Update the "time_taken" timestamps on the data sheet (rows 2..61) ---
$newTimes = @(
    "2021-10-05 14:22:23.950983",
    "2021-10-05 14:22:23.950990",
    "2021-10-05 14:22:23.950993",
    "2021-10-05 14:22:23.950995",
    "2021-10-05 14:22:23.950997",
    "2021-10-05 14:22:23.950999",
    "2021-10-05 14:22:23.951001",
    "2021-10-05 14:22:23.951003",
    "2021-10-05 14:22:23.951005",
    "2021-10-05 14:22:23.951007",
    "2021-10-05 14:22:23.951009",
    "2021-10-05 14:22:23.951011",
    "2021-10-05 14:22:23.951013",
    "2021-10-05 14:22:23.951016",
    "2021-10-05 14:22:23.951018",
    "2021-10-05 14:22:23.951020",
    "2021-10-05 14:22:23.951022",
    "2021-10-05 14:22:23.951024",
    "2021-10-05 14:22:23.951026",
    "2021-10-05 14:22:23.951028",
    "2021-10-05 14:22:23.951030",
    "2021-10-05 14:22:23.951032",
    "2021-10-05 14:22:23.951034",
    "2021-10-05 14:22:23.951036",
    "2021-10-05 14:22:23.951039",
    "2021-10-05 14:22:23.951041",
    "2021-10-05 14:22:23.951043",
    "2021-10-05 14:22:23.951045",
    "2021-10-05 14:22:23.951047",
    "2021-10-05 14:22:23.951049",
    "2021-10-05 14:22:23.951051",
    "2021-10-05 14:22:23.951053",
    "2021-10-05 14:22:23.951056",
    "2021-10-05 14:22:23.951058",
    "2021-10-05 14:22:23.951060",
    "2021-10-05 14:22:23.951062",
    "2021-10-05 14:22:23.951064",
    "2021-10-05 14:22:23.951066",
    "2021-10-05 14:22:23.951068",
    "2021-10-05 14:22:23.951070",
    "2021-10-05 14:22:23.951072",
    "2021-10-05 14:22:23.951074",
    "2021-10-05 14:22:23.951076",
    "2021-10-05 14:22:23.951078",
    "2021-10-05 14:22:23.951080",
    "2021-10-05 14:22:23.951082",
    "2021-10-05 14:22:23.951084",
    "2021-10-05 14:22:23.951086",
    "2021-10-05 14:22:23.951088",
    "2021-10-05 14:22:23.951090",
    "2021-10-05 14:22:23.951092",
    "2021-10-05 14:22:23.951095",
    "2021-10-05 14:22:23.951097",
    "2021-10-05 14:22:23.951099",
    "2021-10-05 14:22:23.951101",
    "2021-10-05 14:22:23.951104",
    "2021-10-05 14:22:23.951106",
    "2021-10-05 14:22:23.951108",
    "2021-10-05 14:22:23.951110",
    "2021-10-05 14:22:23.951112"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2. Add the new "metadata" worksheet after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Copy the header cell formatting (bold, border, centered) from the data
# sheet's header row so the new sheet reuses the same style.
$dataSheet.Range("B1").Copy($metaSheet.Range("B1"))
$dataSheet.Range("B1").Copy($metaSheet.Range("C1"))
$dataSheet.Range("B1").Copy($metaSheet.Range("D1"))
$dataSheet.Range("B1").Copy($metaSheet.Range("E1"))
$dataSheet.Range("B1").Copy($metaSheet.Range("F1"))
$dataSheet.Range("B1").Copy($metaSheet.Range("G1"))

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Index cell (A2) mirrors the "data" sheet's index-column style.
$dataSheet.Range("A2").Copy($metaSheet.Range("A2"))
$metaSheet.Range("A2").Value = 0

# Data row
$metaSheet.Range("B2").Value = "Radial dysplasia"
$metaSheet.Range("C2").Value = 247

# "1.15" looks numeric, so Excel would silently coerce it to a Number;
# force Text formatting first (then drop the now-unneeded numFmt) so it
# round-trips as the literal string "1.15", matching the source panel
# version field.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.15"
$metaSheet.Range("D2").ClearFormats()

$metaSheet.Range("E2").Value = "2021-08-17T13:20:42.763498Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:22:23.948485"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/247/?format=json"

# Re-select the data sheet as the active tab (matches original workbook view).
$dataSheet.Activate()
$dataSheet.Range("A1").Select()
